$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column D: "Correct_answer" header, "l" for Purple rows (2-41), "s" for Blue rows (42-101)
$ws.Range("D1").Value = "Correct_answer"

for ($r = 2; $r -le 41; $r++) {
    $ws.Cells.Item($r, 4).Value = "l"
}

for ($r = 42; $r -le 101; $r++) {
    $ws.Cells.Item($r, 4).Value = "s"
}

# Update view/selection state to match the author's final position in the sheet
$ws.Range("D42:D101").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 82
$win.ScrollColumn = 1
